# Update the "想去人数" (want-to-go count) column F values in the
# "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) worksheets to
# reflect refreshed counts as of the regenerated gh-pages data dump
# (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 467
$ws1.Range("F4").Value  = 7824
$ws1.Range("F5").Value  = 92
$ws1.Range("F6").Value  = 208
$ws1.Range("F8").Value  = 29
$ws1.Range("F10").Value = 453
$ws1.Range("F11").Value = 164
$ws1.Range("F13").Value = 441
$ws1.Range("F14").Value = 66
$ws1.Range("F15").Value = 67
$ws1.Range("F16").Value = 25
$ws1.Range("F17").Value = 5692
$ws1.Range("F18").Value = 163
$ws1.Range("F19").Value = 231
$ws1.Range("F20").Value = 1309
$ws1.Range("F22").Value = 347

# --- Sheet "演出" -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 3

# --- Sheet "全部类型" ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 467
$ws4.Range("F4").Value  = 7824
$ws4.Range("F5").Value  = 92
$ws4.Range("F6").Value  = 208
$ws4.Range("F8").Value  = 29
$ws4.Range("F10").Value = 453
$ws4.Range("F11").Value = 164
$ws4.Range("F13").Value = 441
$ws4.Range("F14").Value = 66
$ws4.Range("F15").Value = 67
$ws4.Range("F16").Value = 25
$ws4.Range("F17").Value = 3
$ws4.Range("F18").Value = 5692
$ws4.Range("F20").Value = 163
$ws4.Range("F21").Value = 231
$ws4.Range("F22").Value = 1309
$ws4.Range("F24").Value = 347
